$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force text storage so numeric-looking strings (e.g. "1.01") are not
    # reinterpreted as numbers - matches the original inlineStr cell type.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "42.922.10"
Set-TextValue "E2" "  +1.45%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.291.27"
Set-TextValue "E3" "  -0.50%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.01"
Set-TextValue "E4" "  +0.46%  "

# Row 5 - BNB
Set-TextValue "E5" "  -0.59%  "

# Row 6 - Solana
Set-TextValue "D6" "105.04"
Set-TextValue "E6" "  +1.54%  "

# Row 7 - XRP
Set-TextValue "E7" "  +0.52%  "

# Row 8 - USDC
Set-TextValue "E8" "  +0.16%  "

# Row 9 - Cardano
Set-TextValue "E9" "  -0.23%  "

# Row 10 - Avalanche
Set-TextValue "D10" "39.75"
Set-TextValue "E10" "  +0.16%  "

# Row 11 - Dogecoin
Set-TextValue "E11" "  -0.42%  "

# Row 12 - Polkadot
Set-TextValue "D12" "8.43"

# Row 13 - TRON
Set-TextValue "E13" "  +2.28%  "

# Row 14 - Polygon
Set-TextValue "D14" "1.00"

# Row 15 - Chainlink
Set-TextValue "D15" "15.26"
Set-TextValue "E15" "  +0.07%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "2.638.89"
Set-TextValue "E16" "  -0.49%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.331.05"
Set-TextValue "E17" "  +1.58%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "42.797.58"
Set-TextValue "E18" "  +0.88%  "

# Row 19 - Uniswap
Set-TextValue "D19" "7.42"
Set-TextValue "E19" "  -0.67%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextValue "D20" "13.87"
Set-TextValue "E20" "  +21.96%  "

# Row 21 - ShibaInu
Set-TextValue "E21" "  -0.21%  "

# Row 22 - Litecoin
Set-TextValue "D22" "73.92"
Set-TextValue "E22" "  +0.59%  "

# Row 23 - PancakeSwap
Set-TextValue "D23" "3.56"
Set-TextValue "E23" "  +0.87%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "265.96"
Set-TextValue "E24" "  -3.89%  "

# Row 25 - ImmutableX
Set-TextValue "E25" "  -2.78%  "

# Row 26 - Dai
Set-TextValue "E26" "  +0.32%  "

# Row 27 - now Cosmos (was Filecoin)
Set-TextValue "B27" "Cosmos"
Set-TextValue "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "10.87"
Set-TextValue "E27" "  +0.21%  "

# Row 28 - now Filecoin (was Cosmos)
Set-TextValue "B28" "Filecoin"
Set-TextValue "C28" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D28" "7.30"
Set-TextValue "E28" "  +24.22%  "

# Row 29 - Toncoin
Set-TextValue "E29" "  -0.31%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "22.53"
Set-TextValue "E30" "  -0.93%  "

# Row 31 - InjectiveProtocol
Set-TextValue "D31" "37.22"
Set-TextValue "E31" "  +0.18%  "

# Row 32 - Monero
Set-TextValue "D32" "167.28"
Set-TextValue "E32" "  +1.10%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0876"
Set-TextValue "E33" "  +0.19%  "

# Row 34 - Stellar
Set-TextValue "E34" "  -2.79%  "

# Row 35 - WEMIXToken
Set-TextValue "E35" "  -0.50%  "

# Row 36 - Kaspa
Set-TextValue "E36" "  -4.32%  "

# Row 37 - RenderToken
Set-TextValue "E37" "  -0.51%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.0351"
Set-TextValue "E38" "  -4.19%  "

# Row 39 - NEARProtocol
Set-TextValue "D39" "3.81"
Set-TextValue "E39" "  +2.75%  "

# Row 40 - LidoDAOToken
Set-TextValue "E40" "  -3.61%  "

# Row 41 - ARBITRUM
Set-TextValue "E41" "  +4.55%  "

# Row 42 - MultiversX
Set-TextValue "D42" "70.74"
Set-TextValue "E42" "  +1.14%  "

# Row 43 - Algorand
Set-TextValue "E43" "  +2.53%  "

# Row 44 - BitcoinSV
Set-TextValue "D44" "94.39"

# Row 45 - FirstDigitalUSD
Set-TextValue "E45" "  +0.13%  "

# Row 46 - Celestia
Set-TextValue "D46" "12.19"
Set-TextValue "E46" "  +0.95%  "

# Row 47 - Maker
Set-TextValue "D47" "1.741.27"
Set-TextValue "E47" "  +9.55%  "

# Row 48 - Aave
Set-TextValue "D48" "113.70"
Set-TextValue "E48" "  +0.66%  "

# Row 49 - ordi
Set-TextValue "D49" "80.11"
Set-TextValue "E49" "  -1.40%  "

# Row 50 - FraxShare
Set-TextValue "D50" "8.76"

# Row 51 - THORChain
Set-TextValue "D51" "5.20"
Set-TextValue "E51" "  -0.65%  "
